# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets
# to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 2125   # was 2122
$wsExhibit.Range("F5").Value = 1363   # was 1350
$wsExhibit.Range("F6").Value = 369    # was 366

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 2125   # was 2122
$wsAll.Range("F7").Value = 1363   # was 1350
$wsAll.Range("F8").Value = 369    # was 366
